$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new course column: header + rating value
$ws.Range("I1").Value = "Introduction to Python"
$ws.Range("I2").Value = 5

# Match the other header columns, which are sized to fit their text
# ("best fit" width, same as the rest of the header row).
$ws.Columns.Item(9).ColumnWidth = 20.6
